$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.272399544715881
$ws.Range("B1").Value = 2.312768936157227
$ws.Range("C1").Value = 3.852398872375488
$ws.Range("D1").Value = 2.807035684585571
$ws.Range("E1").Value = 1.347347021102905
